# Chapter 13/14 position_data.xlsx update:
# - rename the dated post/post-code labels from 03/05/2022 to 03/11/2022
# - move the active-cell selection down one row (A3 -> A4)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Post03112022"
$ws.Range("A2").Value = "PostCode03112022"

$ws.Range("A4").Select()
